# The "Northeast Atlantic" (rows 2-28) and "Central North Atlantic" (rows 29-55)
# blocks of data got swapped by mistake; this restores the correct pairing
# between the Area label and its D:K statistics by swapping the two blocks
# back into place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstBlockStart = 2
$lastBlockStart  = 29
$blockSize       = 27   # rows 2-28 and rows 29-55

$cols = @("D", "E", "F", "G", "H", "I", "J", "K")

for ($i = 0; $i -lt $blockSize; $i++) {
    $r1 = $firstBlockStart + $i
    $r2 = $lastBlockStart + $i

    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"

        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2

        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }

    # Swap the Area text labels (column A) to match the now-swapped stats
    $ws.Range("A$r1").Value2 = "Central North Atlantic"
    $ws.Range("A$r2").Value2 = "Northeast Atlantic"
}
